$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Four pairs of rows had their match data (columns F..V) swapped between
#    them (the "home" row got the other's data, and vice-versa). Columns
#    A..E (index / pais / torneio / temporada / data_partida) are untouched.
# ---------------------------------------------------------------------------
function Swap-Rows($r1, $r2) {
    for ($c = 6; $c -le 22; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value2
        $v2 = $ws.Cells.Item($r2, $c).Value2
        $ws.Cells.Item($r1, $c).Value = $v2
        $ws.Cells.Item($r2, $c).Value = $v1
    }
}

Swap-Rows 31 32
Swap-Rows 36 37
Swap-Rows 40 41
Swap-Rows 56 57

# ---------------------------------------------------------------------------
# 2) Eight new match rows were appended at the end of the sheet (81..88),
#    extending the used range from A1:V80 to A1:V88. Copy the formatting of
#    the last existing row (80) down first so the new rows pick up the same
#    cell styles (bold/border/centered index column, datetime number format
#    on the match-date column), then fill in the values.
# ---------------------------------------------------------------------------
$ws.Range("A80:V80").Copy()
$ws.Range("A81:V88").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$newRows = @(
    @{ Row=81; A=80; E=45226.83333333334;  F="Vitesse";          G=1; H="Zwolle";           I=1; J=2.25; K="22/10/2023 14:42"; L=2.19; M="27/10/2023 19:51"; N=3.46; O="22/10/2023 14:42"; P=3.68; Q="27/10/2023 19:55"; R=3.27; S="22/10/2023 14:42"; T=3.33; U="27/10/2023 19:55"; V="https://www.betexplorer.com/football/netherlands/eredivisie/vitesse-zwolle/rTPDy9th/" },
    @{ Row=82; A=81; E=45227.78125;        F="Heerenveen";       G=3; H="Heracles";         I=0; J=1.98; K="22/10/2023 14:42"; L=1.74; M="28/10/2023 18:41"; N=3.91; O="22/10/2023 14:42"; P=4.15; Q="28/10/2023 18:41"; R=3.56; S="22/10/2023 14:42"; T=4.66; U="28/10/2023 18:41"; V="https://www.betexplorer.com/football/netherlands/eredivisie/heerenveen-heracles/GtOLZnB4/" },
    @{ Row=83; A=82; E=45227.83333333334;  F="Almere City";      G=0; H="G.A. Eagles";      I=0; J=2.51; K="22/10/2023 17:12"; L=3.13; M="28/10/2023 19:58"; N=3.48; O="22/10/2023 17:12"; P=3.57; Q="28/10/2023 19:58"; R=2.89; S="22/10/2023 17:12"; T=2.34; U="28/10/2023 19:58"; V="https://www.betexplorer.com/football/netherlands/eredivisie/almere-city-g-a-eagles/xjPHzTeb/" },
    @{ Row=84; A=83; E=45227.83333333334;  F="Sparta Rotterdam";  G=2; H="Waalwijk";         I=0; J=1.75; K="22/10/2023 17:12"; L=1.78; M="28/10/2023 19:35"; N=3.96; O="22/10/2023 17:12"; P=3.95; Q="28/10/2023 19:37"; R=4.62; S="22/10/2023 17:12"; T=4.63; U="28/10/2023 19:35"; V="https://www.betexplorer.com/football/netherlands/eredivisie/sparta-rotterdam-waalwijk/YHIQY6QA/" },
    @{ Row=85; A=84; E=45227.875;          F="Sittard";          G=0; H="Utrecht";          I=0; J=2.41; K="22/10/2023 12:43"; L=2.38; M="28/10/2023 20:56"; N=3.74; O="22/10/2023 12:43"; P=3.52; Q="28/10/2023 20:56"; R=2.8;  S="22/10/2023 12:43"; T=3.1;  U="28/10/2023 20:52"; V="https://www.betexplorer.com/football/netherlands/eredivisie/sittard-utrecht/C8HUXQuH/" },
    @{ Row=86; A=85; E=45228.51041666666;  F="Twente";           G=2; H="Feyenoord";        I=1; J=2.62; K="22/10/2023 14:42"; L=3.66; M="29/10/2023 12:13"; N=3.67; O="22/10/2023 14:42"; P=3.79; Q="29/10/2023 12:13"; R=2.6;  S="22/10/2023 14:42"; T=2.04; U="29/10/2023 12:12"; V="https://www.betexplorer.com/football/netherlands/eredivisie/twente-feyenoord/6VGYWpeN/" },
    @{ Row=87; A=86; E=45228.60416666666;  F="PSV";               G=5; H="Ajax";             I=2; J=1.54; K="22/10/2023 12:43"; L=1.29; M="29/10/2023 14:23"; N=4.9;  O="22/10/2023 12:43"; P=6.33; Q="29/10/2023 14:29"; R=5.21; S="22/10/2023 12:43"; T=9.79; U="29/10/2023 14:29"; V="https://www.betexplorer.com/football/netherlands/eredivisie/psv-ajax/I5KxW4AT/" },
    @{ Row=88; A=87; E=45228.60416666666;  F="FC Volendam";      G=3; H="Excelsior";        I=1; J=3.04; K="22/10/2023 14:42"; L=2.62; M="29/10/2023 14:26"; N=3.88; O="22/10/2023 14:42"; P=3.81; Q="29/10/2023 14:24"; R=2.25; S="22/10/2023 14:42"; T=2.61; U="29/10/2023 14:21"; V="https://www.betexplorer.com/football/netherlands/eredivisie/fc-volendam-excelsior/Onc5ROIp/" }
)

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value  = $row.A
    $ws.Cells.Item($r, 2).Value  = "netherlands"
    $ws.Cells.Item($r, 3).Value  = "eredivisie"
    $ws.Cells.Item($r, 4).Value  = "2023-2024"
    $ws.Cells.Item($r, 5).Value  = $row.E
    $ws.Cells.Item($r, 6).Value  = $row.F
    $ws.Cells.Item($r, 7).Value  = $row.G
    $ws.Cells.Item($r, 8).Value  = $row.H
    $ws.Cells.Item($r, 9).Value  = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $ws.Cells.Item($r, 21).Value = $row.U
    $ws.Cells.Item($r, 22).Value = $row.V
}
